$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = "Thu Jan 25 17:22:51 EST 2024"
$ws.Cells.Item(3, 2).Value = "Thu Jan 25 17:23:01 EST 2024"
$ws.Cells.Item(4, 2).Value = "Thu Jan 25 17:23:11 EST 2024"
$ws.Cells.Item(5, 2).Value = "Thu Jan 25 17:23:21 EST 2024"
$ws.Cells.Item(6, 2).Value = "Thu Jan 25 17:23:30 EST 2024"
$ws.Cells.Item(7, 2).Value = "Thu Jan 25 17:23:40 EST 2024"
$ws.Cells.Item(8, 2).Value = "Thu Jan 25 17:23:50 EST 2024"
$ws.Cells.Item(9, 2).Value = "Thu Jan 25 17:24:00 EST 2024"
$ws.Cells.Item(10, 2).Value = "Thu Jan 25 17:24:09 EST 2024"
$ws.Cells.Item(11, 2).Value = "Thu Jan 25 17:24:19 EST 2024"
$ws.Cells.Item(12, 2).Value = "Thu Jan 25 17:24:29 EST 2024"
$ws.Cells.Item(13, 2).Value = "Thu Jan 25 17:24:39 EST 2024"
$ws.Cells.Item(14, 2).Value = "Thu Jan 25 17:24:49 EST 2024"
$ws.Cells.Item(15, 2).Value = "Thu Jan 25 17:25:01 EST 2024"
$ws.Cells.Item(16, 2).Value = "Thu Jan 25 17:25:12 EST 2024"
$ws.Cells.Item(17, 2).Value = "Thu Jan 25 17:25:22 EST 2024"
$ws.Cells.Item(18, 2).Value = "Thu Jan 25 17:25:32 EST 2024"
$ws.Cells.Item(25, 2).Value = "Thu Jan 25 17:25:43 EST 2024"
$ws.Cells.Item(26, 2).Value = "Thu Jan 25 17:25:53 EST 2024"
$ws.Cells.Item(27, 2).Value = "Thu Jan 25 17:26:03 EST 2024"
$ws.Cells.Item(28, 2).Value = "Thu Jan 25 17:26:12 EST 2024"
$ws.Cells.Item(29, 2).Value = "Thu Jan 25 17:26:23 EST 2024"
$ws.Cells.Item(30, 2).Value = "Thu Jan 25 17:26:33 EST 2024"
$ws.Cells.Item(31, 2).Value = "Thu Jan 25 17:26:43 EST 2024"
$ws.Cells.Item(32, 2).Value = "Thu Jan 25 17:26:53 EST 2024"
$ws.Cells.Item(33, 2).Value = "Thu Jan 25 17:27:02 EST 2024"
$ws.Cells.Item(34, 2).Value = "Thu Jan 25 17:27:12 EST 2024"
$ws.Cells.Item(35, 2).Value = "Thu Jan 25 17:27:22 EST 2024"
$ws.Cells.Item(36, 2).Value = "Thu Jan 25 17:27:32 EST 2024"
$ws.Cells.Item(37, 2).Value = "Thu Jan 25 17:27:41 EST 2024"
$ws.Cells.Item(38, 2).Value = "Thu Jan 25 17:27:51 EST 2024"
$ws.Cells.Item(39, 2).Value = "Thu Jan 25 17:28:01 EST 2024"
$ws.Cells.Item(40, 2).Value = "Thu Jan 25 17:28:10 EST 2024"
$ws.Cells.Item(41, 2).Value = "Thu Jan 25 17:28:20 EST 2024"
$ws.Cells.Item(42, 2).Value = "Thu Jan 25 17:28:30 EST 2024"
$ws.Cells.Item(43, 2).Value = "Thu Jan 25 17:28:40 EST 2024"
$ws.Cells.Item(44, 2).Value = "Thu Jan 25 17:28:50 EST 2024"
$ws.Cells.Item(45, 2).Value = "Thu Jan 25 17:29:00 EST 2024"
$ws.Cells.Item(46, 2).Value = "Thu Jan 25 17:29:09 EST 2024"
$ws.Cells.Item(47, 2).Value = "Thu Jan 25 17:29:19 EST 2024"
$ws.Cells.Item(48, 2).Value = "Thu Jan 25 17:29:29 EST 2024"
